$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12:D12").Copy()
$ws.Range("A13:D13").PasteSpecial(-4122)

$ws.Range("A13").Value = "Problema ao cadastrar um radical para uma patente"
$ws.Range("B13").Value = "Defeito"
$ws.Range("C13").Value = "Em análise"
$ws.Range("D13").Value = ""

$ws.Range("A13:D13").Select()
